$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numeric-looking strings (e.g. "26.704.53")
# that must stay as plain text, exactly like the original inline strings.
# Temporarily force a text number format on the column while assigning the
# new values, then clear the formats again so cells end up with no style
# (matching the original, unstyled D-column cells).
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.704.53"
$ws.Range("E2").Value = "  +0.39%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.599.48"
$ws.Range("E3").Value = "  +1.00%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "211.19"
$ws.Range("E5").Value = "  +0.25%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.514"
$ws.Range("E6").Value = "  +1.47%  "

# Row 8 - Dogecoin
$ws.Range("D8").Value = "0.0617"
$ws.Range("E8").Value = "  -0.03%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -1.20%  "

# Row 10 - Solana
$ws.Range("E10").Value = "  +0.42%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +1.08%  "

# Rows 12 / 13 - Wrapped Ether pair swapped content, with updated price/volume
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.675.46"
$ws.Range("E12").Value = "  +5.82%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.823.53"
$ws.Range("E13").Value = "  +1.01%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -0.59%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "64.97"
$ws.Range("E16").Value = "  +0.71%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "26.682.76"
$ws.Range("E17").Value = "  +0.25%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "0.0₃0728"
$ws.Range("E18").Value = "  -0.18%  "

# Row 19
$ws.Range("D19").Value = "208.46"
$ws.Range("E19").Value = "  +0.35%  "

# Row 20
$ws.Range("E20").Value = "  +0.01%  "

# Row 21
$ws.Range("D21").Value = "6.79"
$ws.Range("E21").Value = "  +0.45%  "

# Row 22
$ws.Range("E22").Value = "  +0.11%  "

# Row 23
$ws.Range("E23").Value = "  -2.34%  "

# Row 24
$ws.Range("D24").Value = "8.86"
$ws.Range("E24").Value = "  -0.31%  "

# Row 25
$ws.Range("D25").Value = "145.92"
$ws.Range("E25").Value = "  -0.27%  "

# Row 26
$ws.Range("E26").Value = "  +0.35%  "

# Row 27
$ws.Range("E27").Value = "  -2.42%  "

# Row 28
$ws.Range("E28").Value = "  +1.88%  "

# Row 29
$ws.Range("D29").Value = "15.26"
$ws.Range("E29").Value = "  -0.14%  "

# Row 30
$ws.Range("E30").Value = "  +0.72%  "

# Row 31
$ws.Range("E31").Value = "  +0.16%  "

# Row 32
$ws.Range("D32").Value = "3.23"
$ws.Range("E32").Value = "  -0.67%  "

# Row 33
$ws.Range("D33").Value = "0.661"
$ws.Range("E33").Value = "  -2.38%  "

# Row 34
$ws.Range("E34").Value = "  +0.58%  "

# Row 35
$ws.Range("D35").Value = "1.289.14"
$ws.Range("E35").Value = "  -2.26%  "

# Row 36
$ws.Range("E36").Value = "  -2.03%  "

# Row 37
$ws.Range("E37").Value = "  -0.71%  "

# Row 38
$ws.Range("E38").Value = "  -0.52%  "

# Row 39
$ws.Range("E39").Value = "  +2.30%  "

# Row 40
$ws.Range("E40").Value = "  +0.03%  "

# Row 41
$ws.Range("D41").Value = "5.42"
$ws.Range("E41").Value = "  +1.52%  "

# Row 42
$ws.Range("E42").Value = "  +0.95%  "

# Row 43
$ws.Range("E43").Value = "  +0.25%  "

# Row 44
$ws.Range("D44").Value = "63.56"
$ws.Range("E44").Value = "  +0.24%  "

# Row 45
$ws.Range("D45").Value = "1.736.07"
$ws.Range("E45").Value = "  +1.01%  "

# Row 46
$ws.Range("D46").Value = "0.903"
$ws.Range("E46").Value = "  +8.93%  "

# Row 47
$ws.Range("D47").Value = "90.00"
$ws.Range("E47").Value = "  +0.66%  "

# Row 48
$ws.Range("E48").Value = "  -0.60%  "

# Row 49
$ws.Range("E49").Value = "  +2.24%  "

# Row 50
$ws.Range("D50").Value = "0.0505"
$ws.Range("E50").Value = "  -0.22%  "

# Row 51
$ws.Range("D51").Value = "7.46"
$ws.Range("E51").Value = "  +0.14%  "

# Restore the D column cells to their original unstyled state.
$priceCol.ClearFormats()
